$d = $word.ActiveDocument

$replacements = @(
    @("571÷2=", "838÷9="),
    @("385÷6=", "654÷8="),
    @("544÷5=", "144÷7="),
    @("516÷9=", "514÷7="),
    @("702÷8=", "102÷6="),
    @("574÷4=", "141÷3="),
    @("350÷8=", "987÷6="),
    @("993÷2=", "124÷4="),
    @("900÷4=", "417÷9="),
    @("578÷8=", "458÷3="),
    @("429÷2=", "677÷7="),
    @("874÷2=", "732÷2="),
    @("142÷8=", "233÷3="),
    @("310÷6=", "422÷7="),
    @("646÷6=", "475÷9="),
    @("509÷4=", "651÷7="),
    @("196÷4=", "133÷9="),
    @("555÷2=", "689÷2="),
    @("471÷8=", "797÷8="),
    @("542÷8=", "961÷5="),
    @("648÷2=", "258÷9="),
    @("890÷3=", "430÷8="),
    @("128÷6=", "450÷8="),
    @("732÷7=", "796÷5="),
    @("588÷7=", "285÷5=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
